$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.261.39"
$ws.Range("E2").Value = "  +1.89%  "
$ws.Range("D3").Value = "3.589.77"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "567.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.610"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.69%  "
$ws.Range("E8").Value = "  -0.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.679"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "63.74"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +14.08%  "
$ws.Range("E11").Value = "  -0.52%  "
$ws.Range("E12").Value = "  +4.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.07"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.41%  "
$ws.Range("D14").Value = "4.157.77"
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("D15").Value = "3.603.71"
$ws.Range("E15").Value = "  +0.60%  "
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.126"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.13"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.78%  "
$ws.Range("D18").Value = "67.978.14"
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "402.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.62"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "736.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +15.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.63"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "31.48"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.09"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "63.59"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.53%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.113"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "41.54"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.421"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.96%  "
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.29"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +11.51%  "
$ws.Range("E39").Value = "  +28.78%  "
$ws.Range("D40").Value = "0.0₃0743"
$ws.Range("E40").Value = "  -2.32%  "
$ws.Range("D41").Value = "3.185.79"
$ws.Range("E41").Value = "  -1.40%  "
$ws.Range("E42").Value = "  -1.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.997"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.62"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0411"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.130"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.09"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "139.81"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.69"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.40%  "
